# Applies the "fixed workflow" re-run of the day/night sensitivity study:
# the first 4 data points (Cutoff = 0..3) are dropped from each results
# sheet (NBR and BAR) and the remaining rows shift up, while the Cutoff
# index column (A) is renumbered back to a contiguous 0-based sequence.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    # Remove the 4 rows right under the header (Excel rows 2-5, which hold
    # Cutoff values 0,1,2,3). This shifts every row below them up by 4 and
    # carries the B/C values along with them.
    $ws.Range("A2:A5").EntireRow.Delete() | Out-Null

    # After the shift, column A (Cutoff) must be a fresh contiguous index
    # starting again at 0, matching the target data.
    $lastRow = $ws.UsedRange.Rows.Count
    for ($r = 2; $r -le $lastRow; $r++) {
        $ws.Cells.Item($r, 1).Value = $r - 2
    }
}
